$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: a teacher/admin account row appended below the existing
# student rows (nisn, fullName, username, password, classGroup).
$ws.Range("A4").Value = 654321
$ws.Range("B4").Value = "Ahmad Zidan"
$ws.Range("C4").Value = "zidan"
$ws.Range("D4").Value = "Zidan456"
$ws.Range("E4").Value = "ADMIN"

# Match the formatting already used by the rest of the table (B1:E3) so the
# new row's text cells share the same style instead of minting a new one.
$ws.Range("B3:E3").Copy()
$ws.Range("B4:E4").PasteSpecial(-4122)

# Leave the selection on the last cell that was filled in, like a user who
# just finished typing the new row.
[void]$ws.Range("E4").Select()
